$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths
# (COM's ColumnWidth round-trips ~5/6 of a character wider than the stored
# OOXML "width" attribute, so compensate to land on the exact target width)
$ws.Columns.Item(2).ColumnWidth = 27 - 5/6
$ws.Columns.Item(3).ColumnWidth = 32 - 5/6
$ws.Columns.Item(5).ColumnWidth = 13 - 5/6

# Add new rows of enrollment data
# Phone numbers keep a leading zero, so format those cells as text first
$ws.Range("E6:E7").NumberFormat = "@"

$ws.Range("A6").Value = "2025-05-21 15:29:45"
$ws.Range("B6").Value = "LOKESWAR RAJU GUNDLAPALLI"
$ws.Range("C6").Value = "lokeshgundlapalli143@gmail.com"
$ws.Range("D6").Value = "MLOps"
$ws.Range("E6").Value = "08374705188"

$ws.Range("A7").Value = "2025-05-21 15:30:12"
$ws.Range("B7").Value = "LOKESWAR RAJU GUNDLAPALLI"
$ws.Range("C7").Value = "lokeshgundlapalli143@gmail.com"
$ws.Range("D7").Value = "Machine Learning"
$ws.Range("E7").Value = "08374705188"
